$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (border/bold/alignment) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 5

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1

$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 4

$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 3
